# Refresh cryptos list values/percentages scraped by the Actions job.
# Source data is plain text (coin names, links, price strings, volume
# percentages) so every cell below is written as text. Column D prices
# are apostrophe-prefixed (Excel's own "store as text" convention) so
# purely-numeric-looking prices like "573.22" are not coerced to numbers,
# then the style is reset to Normal so the quote-prefix flag left behind
# by that trick does not change the cell's formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.509.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.57%  "

$ws.Range("D3").Value = "'2.480.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'573.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.91%  "

$ws.Range("D6").Value = "'149.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.94%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  +1.88%  "

$ws.Range("E9").Value = "  +4.77%  "

$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("D11").Value = "'0.365"
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = "  +2.98%  "

$ws.Range("E13").Value = "  +6.33%  "

$ws.Range("E14").Value = "  +7.14%  "

$ws.Range("D15").Value = "'2.900.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "

$ws.Range("D16").Value = "'63.580.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.66%  "

$ws.Range("D17").Value = "'2.487.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.56%  "

$ws.Range("E18").Value = "  +2.68%  "

$ws.Range("D19").Value = "'7.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.43%  "

$ws.Range("E20").Value = "  +3.10%  "

$ws.Range("D21").Value = "'329.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.08%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +9.11%  "

$ws.Range("D24").Value = "'67.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "

$ws.Range("D25").Value = "'643.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +15.40%  "

$ws.Range("E26").Value = "  +13.35%  "

$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").Value = "'2.609.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.95%  "

$ws.Range("E29").Value = "  +10.70%  "

$ws.Range("E30").Value = "  +4.69%  "

$ws.Range("D31").Value = "'0.996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("E32").Value = "  -1.74%  "

$ws.Range("D33").Value = "'1.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.13%  "

$ws.Range("D34").Value = "'5.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.63%  "

$ws.Range("E35").Value = "  +3.85%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.388"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.43%  "

$ws.Range("D38").Value = "'5.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.97%  "

$ws.Range("E39").Value = "  +2.57%  "

$ws.Range("E40").Value = "  +2.40%  "

$ws.Range("D41").Value = "'147.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.11%  "

$ws.Range("E42").Value = "  +18.82%  "

$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").Value = "'152.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "

$ws.Range("E45").Value = "  +4.63%  "

$ws.Range("E46").Value = "  +4.99%  "

$ws.Range("D47").Value = "'21.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.12%  "

$ws.Range("E48").Value = "  +3.58%  "

$ws.Range("E49").Value = "  +5.81%  "

$ws.Range("D50").Value = "'0.0930"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.21%  "

$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'0.747"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.85%  "
